# Updated cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '29.936.49'
    'E2' = '  +0.11%  '
    'D3' = '1.896.14'
    'D4' = '1.000'
    'E4' = '  -0.03%  '
    'D5' = '0.7758'
    'E5' = '  -2.16%  '
    'D6' = '244.82'
    'E6' = '  +0.39%  '
    'E7' = '  -0.02%  '
    'D8' = '0.3143'
    'E8' = '  -0.63%  '
    'D9' = '25.81'
    'E9' = '  +1.66%  '
    'D10' = '0.07272'
    'E10' = '  +1.25%  '
    'D11' = '0.08936'
    'E11' = '  +10.31%  '
    'D12' = '0.7756'
    'E12' = '  +1.07%  '
    'D13' = '5.463'
    'E13' = '  -2.22%  '
    'D14' = '94.90'
    'E14' = '  +2.46%  '
    'D15' = '1.849.79'
    'E15' = '  -1.87%  '
    'D16' = '6.205'
    'E16' = '  +0.54%  '
    'D17' = '29.880.45'
    'E17' = '  -0.07%  '
    'D18' = '14.01'
    'E18' = '  +0.45%  '
    'D19' = '246.96'
    'E19' = '  +0.97%  '
    'D20' = '0.000007915'
    'E20' = '  +1.64%  '
    'D21' = '8.156'
    'E21' = '  -0.99%  '
    'D22' = '0.9999'
    'E22' = '  -0.07%  '
    'D23' = '2.119.29'
    'E23' = '  -0.16%  '
    'D24' = '1.000'
    'D25' = '0.1592'
    'E25' = '  -5.01%  '
    'D26' = '9.564'
    'E26' = '  +0.81%  '
    'D27' = '163.20'
    'E27' = '  -0.53%  '
    'D28' = '18.87'
    'D29' = '2.053'
    'E29' = '  -0.39%  '
    'E30' = '  +1.97%  '
    'D31' = '1.547'
    'E31' = '  -0.13%  '
    'E32' = '  +0.93%  '
    'D33' = '4.127'
    'E33' = '  +0.82%  '
    'D34' = '0.05535'
    'E34' = '  -1.07%  '
    'D35' = '1.251'
    'E35' = '  -2.23%  '
    'D36' = '0.7556'
    'E36' = '  +1.90%  '
    'D37' = '0.9987'
    'E37' = '  +0.22%  '
    'D38' = '2.720'
    'E38' = '  +3.48%  '
    'D39' = '0.01973'
    'E39' = '  +2.13%  '
    'D40' = '2.794'
    'E40' = '  +0.43%  '
    'D41' = '0.4529'
    'E41' = '  +2.46%  '
    'D42' = '74.16'
    'E42' = '  -0.14%  '
    'D43' = '6.080'
    'E43' = '  +2.43%  '
    'D44' = '1.086.15'
    'E44' = '  -6.67%  '
    'D45' = '0.8558'
    'E45' = '  +0.35%  '
    'D47' = '1.899'
    'E47' = '  +0.78%  '
    'D48' = '102.85'
    'E48' = '  -1.73%  '
    'D49' = '7.631'
    'E49' = '  +2.28%  '
    'D50' = '9.883'
    'E50' = '  -0.80%  '
    'D51' = '2.999'
    'E51' = '  -0.59%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    if ($cellRef -match "^D") {
        # Price column holds numeric-looking strings (e.g. "1.000", "29.936.49").
        # Force text so Excel does not coerce them into numbers / strip formatting,
        # then reset the style back to Normal so no stray format is left behind.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
